# Rename "Kupno" sheet's sales-related headers from space-separated to
# underscore-separated names (matches the already-underscored style used
# on the "Sprzedaz" sheet / table), so the generated JPK file validates.
$wb = $excel.ActiveWorkbook

$wsSprzedaz = $wb.Worksheets.Item("Sprzedaz")
$wsKupno = $wb.Worksheets.Item("Kupno")

# Fix duplicated ID value on Sprzedaz (row 2 and row 3 both had ID = 1;
# the second data row should be ID = 2).
$wsSprzedaz.Range("A3").Value = 2

# Rename headers on Kupno sheet (also updates the backing table + list
# columns automatically).
$wsKupno.Range("C1").Value = "NIP_sprzedawcy"
$wsKupno.Range("D1").Value = "Nazwa_sprzedawcy"
$wsKupno.Range("E1").Value = "Nr_faktury"
$wsKupno.Range("F1").Value = "Data_faktury"
$wsKupno.Range("G1").Value = "Kwota_netto"
$wsKupno.Range("H1").Value = "Kwota_podatku"

# Move the active tab / selection from Sprzedaz to Kupno.
$wsSprzedaz.Activate()
$wsSprzedaz.Range("B4").Select()

$wsKupno.Activate()
$wsKupno.Range("K13").Select()
